$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.167
$ws.Range("E2").ClearContents()
$ws.Range("G2").Value = 0.1634561067606838
$ws.Range("H2").Value = 0.1634561067606838
$ws.Range("I2").Value = -0.007732596159881378
$ws.Range("J2").Value = -0.007732596159881378
$ws.Range("K2").Value = -17.791
$ws.Range("L2").Value = -0.000954423603999234
$ws.Range("M2").Value = 495
$ws.Range("N2").Value = 0.05350252380592094
$ws.Range("O2").Value = -27.82305660165252
$ws.Range("P2").Value = 275.7
$ws.Range("Q2").Value = 0.02979928447129778
$ws.Range("R2").Value = -15.49659940419313
$ws.Range("S2").Value = 219.3
$ws.Range("T2").Value = 0.4430303030303031
$ws.Range("U2").Value = 4577.03
$ws.Range("V2").Value = 0.4947124374452814
$ws.Range("W2").Value = 0.3215874377407865
$ws.Range("X2").Value = 0.05237640762358182
$ws.Range("Y2").Value = 0.2692110301172047
$ws.Range("Z2").Value = 1.074427972499405
$ws.Range("AA2").Value = 0.1717292152745801
$ws.Range("AB2").Value = 0.04171469586693065
$ws.Range("AC2").Value = 0.1300145194076495
$ws.Range("AD2").Value = 8650.92
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 8650.92
$ws.Range("AG2").Value = 4073.89
$ws.Range("AH2").Value = 0.4832154934250582
$ws.Range("AI2").Value = 0.3454714704856972
$ws.Range("AJ2").Value = 0.3057147080961054
$ws.Range("AK2").Value = 0.1990771056124254
$ws.Range("AL2").Value = 475.999
$ws.Range("AM2").Value = 475.999
$ws.Range("AN2").Value = 40.54041895121608
$ws.Range("AO2").Value = -0.3028157622179879
$ws.Range("AP2").Value = 19.09128825155818
$ws.Range("AQ2").Value = -0.3028157622179879

# Row 3
$ws.Range("D3").Value = 0.177
$ws.Range("G3").Value = 0.01461245235069886
$ws.Range("H3").Value = 0.01461245235069886
$ws.Range("I3").Value = 0.1473951715374841
$ws.Range("J3").Value = 0.1473951715374841
$ws.Range("K3").Value = 0.909
$ws.Range("L3").Value = 0.1155019059720458
$ws.Range("U3").Value = 0.63
$ws.Range("V3").Value = 0.04772727272727273
$ws.Range("W3").Value = 0.6446808510638299
$ws.Range("X3").Value = 0.0427892584120525
$ws.Range("Y3").Value = 0.6018915926517774
$ws.Range("Z3").Value = 2.387018501668183
$ws.Range("AA3").Value = 0.3518350015165302
$ws.Range("AB3").Value = 0.04062323492004773
$ws.Range("AC3").Value = 0.3112117665964825
$ws.Range("AD3").Value = 1.32
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1.32
$ws.Range("AG3").Value = 0.6900000000000001
$ws.Range("AH3").Value = 0.09090909090909091
$ws.Range("AI3").Value = 0.2688391038696538
$ws.Range("AJ3").Value = 0.04967602591792657
$ws.Range("AK3").Value = 0.1612149532710281
$ws.Range("AL3").Value = 0.199
$ws.Range("AM3").Value = 0.199
$ws.Range("AN3").Value = 1.211009174311927
$ws.Range("AO3").Value = 5.829145728643216
$ws.Range("AP3").Value = 0.6330275229357798
$ws.Range("AQ3").Value = 5.829145728643216

# Row 4
$ws.Range("D4").Value = 0.157
$ws.Range("E4").ClearContents()
$ws.Range("G4").Value = 0.1635189747057592
$ws.Range("H4").Value = 0.1635189747057592
$ws.Range("I4").Value = -0.007798118361804785
$ws.Range("J4").Value = -0.007798118361804785
$ws.Range("K4").Value = -18.7
$ws.Range("L4").Value = -0.00100361192956469
$ws.Range("M4").Value = 495
$ws.Range("N4").Value = 0.05357896673774448
$ws.Range("O4").Value = -26.47058823529412
$ws.Range("P4").Value = 275.7
$ws.Range("Q4").Value = 0.02984186086787102
$ws.Range("R4").Value = -14.74331550802139
$ws.Range("S4").Value = 219.3
$ws.Range("T4").Value = 0.4430303030303031
$ws.Range("U4").Value = 4576.4
$ws.Range("V4").Value = 0.4953510775325532
$ws.Range("W4").Value = -0.001505975582256869
$ws.Range("X4").Value = 0.06196355683511115
$ws.Range("Y4").Value = -0.06346953241736802
$ws.Range("Z4").Value = 1.074178484953303
$ws.Range("AA4").Value = -0.00837657096737
$ws.Range("AB4").Value = 0.04280615681381357
$ws.Range("AC4").Value = -0.05118272778118357
$ws.Range("AD4").Value = 8649.6
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 8649.6
$ws.Range("AG4").Value = 4073.200000000001
$ws.Range("AH4").Value = 0.4835339299989378
$ws.Range("AI4").Value = 0.3454864994408052
$ws.Range("AJ4").Value = 0.3059818658493528
$ws.Range("AK4").Value = 0.199085026100217
$ws.Range("AL4").Value = 475.8
$ws.Range("AM4").Value = 475.8
$ws.Range("AN4").Value = 40.74234573716439
$ws.Range("AO4").Value = -0.305380411937789
$ws.Range("AP4").Value = 19.18605746585022
$ws.Range("AQ4").Value = -0.305380411937789
